# Updates cryptos price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.594.69'
$ws.Range("E2").Value = '  +1.14%  '
$ws.Range("D3").Value = '1.796.09'
$ws.Range("E3").Value = '  -1.10%  '
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.33%  '
$ws.Range("D5").Value = "'225.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.77%  '
$ws.Range("D6").Value = "'0.555"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.93%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").Value = "'33.10"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.36%  '
$ws.Range("D9").Value = "'0.285"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.60%  '
$ws.Range("D10").Value = "'0.0665"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.31%  '
$ws.Range("E11").Value = '  +0.20%  '
$ws.Range("D12").Value = '2.052.90'
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("D13").Value = "'11.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.10%  '
$ws.Range("D14").Value = '1.783.57'
$ws.Range("E14").Value = '  -2.06%  '
$ws.Range("D15").Value = "'0.641"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").Value = '34.507.29'
$ws.Range("E16").Value = '  +0.89%  '
$ws.Range("E17").Value = '  +0.28%  '
$ws.Range("D18").Value = "'69.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").Value = "'255.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.80%  '
$ws.Range("D20").Value = '0.0₃0749'
$ws.Range("E20").Value = '  -0.47%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D22").Value = "'10.47"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.07%  '
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  -3.63%  '
$ws.Range("D25").Value = "'157.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").Value = "'16.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.55%  '
$ws.Range("E27").Value = '  +0.72%  '
$ws.Range("E28").Value = '  -2.89%  '
$ws.Range("D29").Value = "'0.998"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.27%  '
$ws.Range("D30").Value = "'3.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.22%  '
$ws.Range("D31").Value = "'0.0517"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.43%  '
$ws.Range("E32").Value = '  -2.06%  '
$ws.Range("D33").Value = "'3.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.02%  '
$ws.Range("D34").Value = "'1.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.45%  '
$ws.Range("D35").Value = '1.461.03'
$ws.Range("E36").Value = '  -2.11%  '
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("D38").Value = "'0.632"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  +2.40%  '
$ws.Range("D40").Value = "'83.06"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("E41").Value = '  -1.65%  '
$ws.Range("E42").Value = '  -1.07%  '
$ws.Range("D43").Value = "'2.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.30%  '
$ws.Range("D44").Value = "'0.0507"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.78%  '
$ws.Range("D45").Value = "'5.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.67%  '
$ws.Range("D46").Value = '1.952.26'
$ws.Range("E46").Value = '  -0.90%  '
$ws.Range("E47").Value = '  -3.56%  '
$ws.Range("D48").Value = "'12.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.78%  '
$ws.Range("D49").Value = "'0.999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.18%  '
$ws.Range("D50").Value = "'99.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.05%  '
$ws.Range("D51").Value = "'49.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.23%  '

Write-Host "Updated" ([string]85) "cells"
